# New crime data collected - weekly CompStat update for 105th Precinct
# Report period rolls forward one week (Vol 29 No 49 -> No 50; 12/5-12/11/2022 -> 12/12-12/18/2022)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates ----
$ws.Range("A8").Value = "Volume 29   Number  50"
$ws.Range("C9").Value = "Report Covering the Week  12/12/2022  Through  12/18/2022"

# ---- Row 14: Murder ----
# F14/G14 flip from "4 vs 1" to "1 vs 0" (G14 becomes the literal text "0", H14 becomes "***.*")
$ws.Range("F14").Value = 1
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))
$ws.Range("N14").Value = -56.521739130434

# ---- Row 15: Rape ----
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("I15").Value = 13
$ws.Range("J15").Value = 42
$ws.Range("K15").Value = -69.047619047619
$ws.Range("L15").Value = -63.888888888888
$ws.Range("M15").Value = -27.777777777777
$ws.Range("N15").Value = -67.5

# ---- Row 16: Robbery ----
$ws.Range("C14").Copy($ws.Range("C16"))
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -26.666666666666
$ws.Range("I16").Value = 172
$ws.Range("J16").Value = 142
$ws.Range("K16").Value = 21.126760563380
$ws.Range("L16").Value = 1.775147928994
$ws.Range("M16").Value = -48.192771084337
$ws.Range("N16").Value = -83.137254901960

# ---- Row 17: Fel. Assault ----
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -16.666666666666
$ws.Range("I17").Value = 403
$ws.Range("J17").Value = 395
$ws.Range("K17").Value = 2.025316455696
$ws.Range("L17").Value = -1.946472019464
$ws.Range("M17").Value = 50.936329588015
$ws.Range("N17").Value = 0.75

# ---- Row 18: Burglary ----
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 18
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 209
$ws.Range("J18").Value = 157
$ws.Range("K18").Value = 33.121019108280
$ws.Range("L18").Value = 1.951219512195
$ws.Range("M18").Value = -46.272493573264
$ws.Range("N18").Value = -86.822194199243

# ---- Row 19: Gr. Larceny ----
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 47
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = -20.338983050847
$ws.Range("I19").Value = 609
$ws.Range("J19").Value = 482
$ws.Range("K19").Value = 26.348547717842
$ws.Range("L19").Value = 6.282722513089
$ws.Range("M19").Value = 39.359267734553
$ws.Range("N19").Value = 6.282722513089

# ---- Row 20: G.L.A. ----
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 23
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -11.538461538461
$ws.Range("I20").Value = 265
$ws.Range("J20").Value = 159
$ws.Range("K20").Value = 66.666666666666
$ws.Range("L20").Value = -5.693950177935
$ws.Range("M20").Value = -28.954423592493
$ws.Range("N20").Value = -92.131828978622

# ---- Row 21: TOTAL ----
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -68
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 149
$ws.Range("H21").Value = -18.791946308724
$ws.Range("I21").Value = 1681
$ws.Range("J21").Value = 1382
$ws.Range("K21").Value = 21.635311143270
$ws.Range("L21").Value = -0.118835412953
$ws.Range("M21").Value = -8.192244675040
$ws.Range("N21").Value = -76.019971469329

# ---- Row 24: Petit Larceny ----
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = -29.545454545454
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 130
$ws.Range("H24").Value = -23.076923076923
$ws.Range("I24").Value = 1464
$ws.Range("J24").Value = 1066
$ws.Range("K24").Value = 37.335834896810
$ws.Range("L24").Value = 50.617283950617
$ws.Range("M24").Value = 84.615384615384

# ---- Row 25: Misd. Assault ----
$ws.Range("C25").Value = 6
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 11.428571428571
$ws.Range("I25").Value = 626
$ws.Range("J25").Value = 490
$ws.Range("K25").Value = 27.755102040816
$ws.Range("L25").Value = 24.206349206349
$ws.Range("M25").Value = 7.191780821917

# ---- Row 26: UCR Rape* ----
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 1
$ws.Range("F26").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = -85.714285714285
$ws.Range("I26").Value = 29
$ws.Range("J26").Value = 55
$ws.Range("K26").Value = -47.272727272727
$ws.Range("L26").Value = -34.090909090909

# ---- Row 27: Other Sex Crimes ----
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 400
$ws.Range("I27").Value = 61
$ws.Range("K27").Value = 24.489795918367
$ws.Range("L27").Value = 19.607843137254

# ---- Row 28: Shooting Vic. ----
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = -50
$ws.Range("N28").Value = -70.588235294117

# ---- Row 29: Shooting Inc. ----
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("N29").Value = -72.602739726027
